$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped from the dataset (RM 232 and SC 92).
# Row 26 "RM 232" is deleted first; after that shift, the old row 28 "SC 92"
# becomes row 27, so we delete row 27 next.
$ws.Rows("26").Delete()
$ws.Rows("27").Delete()

# Update the "missing data" pattern (values that changed from present -> blank
# or blank -> present) on the resulting 33-row table.
$ws.Range("E5").Value = ""
$ws.Range("E11").Value = -7.9

$ws.Range("D19").Value = -15.5
$ws.Range("E19").Value = ""

$ws.Range("D21").Value = ""

$ws.Range("D23").Value = -13.9

$ws.Range("E25").Value = -7.1

$ws.Range("D27").Value = ""

$ws.Range("F28").Value = 17.44

$ws.Range("E29").Value = ""

$ws.Range("F32").Value = ""

$ws.Range("D33").Value = -14.1
